$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.522.66"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.832.18"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'313.01"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "'0.4294"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("D8").Value = "'0.3668"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.07285"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").Value = "'0.8685"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").Value = "'20.70"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "1.859.57"
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("D13").Value = "'5.410"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "'6.547"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "'0.06940"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "'80.74"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "'0.000008929"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'15.44"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "27.649.93"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "'5.151"
$ws.Range("E22").Value = "  +3.37%  "
$ws.Range("E23").Value = "  +5.36%  "
$ws.Range("D24").Value = "2.098.20"
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").Value = "'1.981"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'154.66"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").Value = "'18.90"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").Value = "'5.144"
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("D29").Value = "'114.19"
$ws.Range("E29").Value = "  -4.96%  "
$ws.Range("D30").Value = "'1.842"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("D31").Value = "'0.08852"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "'0.7562"
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'2.997"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.554"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "'1.137"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "'1.089"
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "'2.799"
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("D41").Value = "'0.5098"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "'0.1668"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "'6.592"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "'8.387"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").Value = "'10.46"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").Value = "'106.36"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").Value = "'0.06507"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "'0.4702"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "'1.623"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").Value = "'64.08"
$ws.Range("E51").Value = "  -0.83%  "
